$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 498
$ws.Range("J2").Value = 498
$ws.Range("L2").Value = 498
$ws.Range("N2").Value = -724
$ws.Range("H32").Value = 9097544
$ws.Range("I32").Value = 10600
$ws.Range("J32").Value = 11116864
$ws.Range("K32").Value = 10600
$ws.Range("L32").Value = 11116864
$ws.Range("M32").Value = -10274
$ws.Range("N32").Value = -11117516
$ws.Range("H87").Value = 19999.87
$ws.Range("J87").Value = 19999.87
$ws.Range("L87").Value = 19999.87
$ws.Range("N87").Value = -22495.87
$ws.Range("H90").Value = 19999.87
$ws.Range("J90").Value = 19999.87
$ws.Range("L90").Value = 59999.61
$ws.Range("N90").Value = -72479.61
$ws.Range("H100").Value = 5864.1055
$ws.Range("I100").Value = 3537.1428
$ws.Range("J100").Value = 7221.5
$ws.Range("K100").Value = 3537.1428
$ws.Range("L100").Value = 7221.5
$ws.Range("M100").Value = -2996.1428
$ws.Range("N100").Value = -8303.5
$ws.Range("H115").Value = 1181.8572
$ws.Range("I115").Value = 1054
$ws.Range("K115").Value = 3162
$ws.Range("M115").Value = -1595
$ws.Range("H138").Value = 6063111.5
$ws.Range("I138").Value = 1147.7693
$ws.Range("J138").Value = 7939434
$ws.Range("K138").Value = 3443.3079
$ws.Range("L138").Value = 23818302
$ws.Range("M138").Value = 1696.6921
$ws.Range("N138").Value = -23828582

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5100.8486
$ws.Range("I61").Value = 4296.522
$ws.Range("K61").Value = 4296.522
$ws.Range("M61").Value = -4084.522
$ws.Range("H74").Value = 1243.4546
$ws.Range("I74").Value = 1301.15
$ws.Range("J74").Value = 666.5
$ws.Range("K74").Value = 1301.15
$ws.Range("L74").Value = 666.5
$ws.Range("M74").Value = -427.1500000000001
$ws.Range("N74").Value = -2414.5
$ws.Range("H77").Value = 1243.4546
$ws.Range("I77").Value = 1301.15
$ws.Range("J77").Value = 666.5
$ws.Range("K77").Value = 6505.75
$ws.Range("L77").Value = 3332.5
$ws.Range("M77").Value = -2137.75
$ws.Range("N77").Value = -12068.5
$ws.Range("H132").Value = 5019.6
$ws.Range("I132").Value = 4793.6763
$ws.Range("K132").Value = 14381.0289
$ws.Range("M132").Value = -11851.0289
$ws.Range("H136").Value = 5100.8486
$ws.Range("I136").Value = 4296.522
$ws.Range("K136").Value = 12889.566
$ws.Range("M136").Value = -10339.566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H86").Value = 4477.48
$ws.Range("I86").Value = 3108.5
$ws.Range("K86").Value = 3108.5
$ws.Range("M86").Value = -1985.5
$ws.Range("H89").Value = 4477.48
$ws.Range("I89").Value = 3108.5
$ws.Range("K89").Value = 15542.5
$ws.Range("M89").Value = -9926.5
$ws.Range("H105").Value = 5222.615
$ws.Range("I105").Value = 3179.4
$ws.Range("K105").Value = 3179.4
$ws.Range("M105").Value = -1432.4
$ws.Range("H131").Value = 31663.084
$ws.Range("J131").Value = 31663.084
$ws.Range("L131").Value = 31663.084
$ws.Range("N131").Value = -41743.084
$ws.Range("H134").Value = 6202.364
$ws.Range("I134").Value = 3378
$ws.Range("J134").Value = 7816.2856
$ws.Range("K134").Value = 10134
$ws.Range("L134").Value = 23448.8568
$ws.Range("M134").Value = -7599
$ws.Range("N134").Value = -28518.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 377.5
$ws.Range("J7").Value = 423.07693
$ws.Range("L7").Value = 423.07693
$ws.Range("N7").Value = -649.0769299999999
$ws.Range("H31").Value = 2967.0417
$ws.Range("I31").Value = 2326.2424
$ws.Range("K31").Value = 2326.2424
$ws.Range("M31").Value = -2031.2424
$ws.Range("H34").Value = 2967.0417
$ws.Range("I34").Value = 2326.2424
$ws.Range("K34").Value = 2326.2424
$ws.Range("M34").Value = -2124.2424
$ws.Range("H94").Value = 7411.1177
$ws.Range("I94").Value = 11279.8
$ws.Range("K94").Value = 11279.8
$ws.Range("M94").Value = -10828.8
$ws.Range("H99").Value = 10692917
$ws.Range("I99").Value = 3490250.2
$ws.Range("J99").Value = 14294251
$ws.Range("K99").Value = 3490250.2
$ws.Range("L99").Value = 14294251
$ws.Range("M99").Value = -3488752.2
$ws.Range("N99").Value = -14297247
$ws.Range("H126").Value = 10692917
$ws.Range("I126").Value = 3490250.2
$ws.Range("J126").Value = 14294251
$ws.Range("K126").Value = 10470750.6
$ws.Range("L126").Value = 42882753
$ws.Range("M126").Value = -10468280.6
$ws.Range("N126").Value = -42887693
$ws.Range("H132").Value = 3755.75
$ws.Range("I132").Value = 3135.258
$ws.Range("K132").Value = 9405.773999999999
$ws.Range("M132").Value = -6875.773999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 2799.8333
$ws.Range("J111").Value = 3259.8
$ws.Range("L111").Value = 9779.400000000001
$ws.Range("N111").Value = -15913.4
$ws.Range("H113").Value = 1235.1154
$ws.Range("I113").Value = 1013.2
$ws.Range("J113").Value = 1373.8125
$ws.Range("K113").Value = 3039.6
$ws.Range("L113").Value = 4121.4375
$ws.Range("M113").Value = -869.6000000000004
$ws.Range("N113").Value = -8461.4375
$ws.Range("H122").Value = 1937.3889
$ws.Range("J122").Value = 1998
$ws.Range("L122").Value = 17982
$ws.Range("N122").Value = -22882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6720.1875
$ws.Range("J113").Value = 12799
$ws.Range("L113").Value = 12799
$ws.Range("N113").Value = -17139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5624
$ws.Range("J7").Value = 6644.375
$ws.Range("L7").Value = 6644.375
$ws.Range("N7").Value = -6868.375
$ws.Range("H46").Value = 7126.4062
$ws.Range("I46").Value = 1831.1818
$ws.Range("J46").Value = 9900.096
$ws.Range("K46").Value = 1831.1818
$ws.Range("L46").Value = 9900.096
$ws.Range("M46").Value = -1643.1818
$ws.Range("N46").Value = -10276.096
$ws.Range("H122").Value = 4855.85
$ws.Range("I122").Value = 4320.533
$ws.Range("K122").Value = 12961.599
$ws.Range("M122").Value = -10511.599
$ws.Range("H126").Value = 5624
$ws.Range("J126").Value = 6644.375
$ws.Range("L126").Value = 19933.125
$ws.Range("N126").Value = -24873.125
$ws.Range("H136").Value = 4943.343
$ws.Range("I136").Value = 2756.7058
$ws.Range("J136").Value = 7008.5
$ws.Range("K136").Value = 8270.117400000001
$ws.Range("L136").Value = 21025.5
$ws.Range("M136").Value = -5720.117400000001
$ws.Range("N136").Value = -26125.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4956
$ws.Range("I122").Value = 2569.7917
$ws.Range("K122").Value = 7709.375100000001
$ws.Range("M122").Value = -5259.375100000001
$ws.Range("H136").Value = 2265.2742
$ws.Range("I136").Value = 1737.0577
$ws.Range("J136").Value = 5012
$ws.Range("K136").Value = 5211.1731
$ws.Range("L136").Value = 15036
$ws.Range("M136").Value = -2661.1731
$ws.Range("N136").Value = -20136
